# "added blank excel file" — wipe the sample schedule rows (2-5) back to a
# blank template: drop the class name / link / time / "no app link" marker
# values (keeping the existing cell styles), drop the hyperlinks that were
# attached to the old Link-or-Code cells, and restore the selection to the
# blank template's first data block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the sample data in rows 2-5 (A: class name, B: link/code, D: time,
# F/H: "No App Link" markers) while leaving the number/text formatting on
# the cells alone (ClearContents preserves the style index, just like the
# target file still carries s="1"/s="2"/s="4" on the now-empty B/D cells).
$ws.Range("A2:A5").ClearContents()
$ws.Range("B2:B5").ClearContents()
$ws.Range("D2:D5").ClearContents()
$ws.Range("F2:F5").ClearContents()
$ws.Range("H2:H5").ClearContents()

# The three sample rows had hyperlinks pointing at zoom.us links; remove
# every hyperlink on the sheet so the <hyperlinks> block disappears.
[void]$ws.Hyperlinks.Delete()

# Restore the view to the top of the blank template, selecting A2:H5.
[void]$ws.Range("A2:H5").Select()
